# Auto-generated Excel COM-interop script to apply the Typhon_Profits.xlsx leve pricing update
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 559.875
$ws.Range("J118").Value = 1200
$ws.Range("L118").Value = 3600
$ws.Range("N118").Value = -6914

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1888.9546
$ws.Range("I61").Value = 1592.1842
$ws.Range("K61").Value = 1592.1842
$ws.Range("M61").Value = -1380.1842
$ws.Range("H74").Value = 32259838
$ws.Range("I74").Value = 47619748
$ws.Range("K74").Value = 47619748
$ws.Range("M74").Value = -47618874
$ws.Range("H77").Value = 32259838
$ws.Range("I77").Value = 47619748
$ws.Range("K77").Value = 238098740
$ws.Range("M77").Value = -238094372
$ws.Range("H102").Value = 1841.0834
$ws.Range("I102").Value = 1826.6364
$ws.Range("K102").Value = 1826.6364
$ws.Range("M102").Value = -204.6364000000001
$ws.Range("H110").Value = 1271.9
$ws.Range("I110").Value = 1272.2858
$ws.Range("K110").Value = 1272.2858
$ws.Range("M110").Value = 772.7141999999999
$ws.Range("H132").Value = 13064.174
$ws.Range("I132").Value = 1997.0322
$ws.Range("J132").Value = 35936.266
$ws.Range("K132").Value = 5991.096600000001
$ws.Range("L132").Value = 107808.798
$ws.Range("M132").Value = -3461.096600000001
$ws.Range("N132").Value = -112868.798
$ws.Range("H136").Value = 1888.9546
$ws.Range("I136").Value = 1592.1842
$ws.Range("K136").Value = 4776.5526
$ws.Range("M136").Value = -2226.5526

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 52520
$ws.Range("J59").Value = 52520
$ws.Range("L59").Value = 52520
$ws.Range("N59").Value = -54214
$ws.Range("H70").Value = 93050
$ws.Range("J70").Value = 93050
$ws.Range("L70").Value = 93050
$ws.Range("N70").Value = -93636
$ws.Range("H73").Value = 93050
$ws.Range("J73").Value = 93050
$ws.Range("L73").Value = 93050
$ws.Range("N73").Value = -95078
$ws.Range("H99").Value = 1347
$ws.Range("I99").Value = 1261.1111
$ws.Range("K99").Value = 1261.1111
$ws.Range("M99").Value = 236.8888999999999
$ws.Range("H134").Value = 4082.5945
$ws.Range("I134").Value = 4213.697
$ws.Range("K134").Value = 12641.091
$ws.Range("M134").Value = -10106.091

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 47.2
$ws.Range("J7").Value = 56.5
$ws.Range("L7").Value = 56.5
$ws.Range("N7").Value = -282.5
$ws.Range("H16").Value = 1190.9166
$ws.Range("I16").Value = 1057.4286
$ws.Range("J16").Value = 1377.8
$ws.Range("K16").Value = 1057.4286
$ws.Range("L16").Value = 1377.8
$ws.Range("M16").Value = -770.4286
$ws.Range("N16").Value = -1951.8
$ws.Range("H22").Value = 177.8
$ws.Range("J22").Value = 212.85715
$ws.Range("L22").Value = 212.85715
$ws.Range("N22").Value = -912.85715
$ws.Range("H31").Value = 3382.1187
$ws.Range("I31").Value = 1556.8286
$ws.Range("J31").Value = 6044
$ws.Range("K31").Value = 1556.8286
$ws.Range("L31").Value = 6044
$ws.Range("M31").Value = -1261.8286
$ws.Range("N31").Value = -6634
$ws.Range("H34").Value = 3382.1187
$ws.Range("I34").Value = 1556.8286
$ws.Range("J34").Value = 6044
$ws.Range("K34").Value = 1556.8286
$ws.Range("L34").Value = 6044
$ws.Range("M34").Value = -1354.8286
$ws.Range("N34").Value = -6448
$ws.Range("H86").Value = 15163212
$ws.Range("I86").Value = 2400
$ws.Range("J86").Value = 18532280
$ws.Range("K86").Value = 2400
$ws.Range("L86").Value = 18532280
$ws.Range("M86").Value = -1277
$ws.Range("N86").Value = -18534526
$ws.Range("H89").Value = 15163212
$ws.Range("I89").Value = 2400
$ws.Range("J89").Value = 18532280
$ws.Range("K89").Value = 12000
$ws.Range("L89").Value = 92661400
$ws.Range("M89").Value = -6384
$ws.Range("N89").Value = -92672632
$ws.Range("H113").Value = 1190.9166
$ws.Range("I113").Value = 1057.4286
$ws.Range("J113").Value = 1377.8
$ws.Range("K113").Value = 1057.4286
$ws.Range("L113").Value = 1377.8
$ws.Range("M113").Value = 1112.5714
$ws.Range("N113").Value = -5717.8
$ws.Range("H134").Value = 1014.0968
$ws.Range("I134").Value = 815.0714
$ws.Range("K134").Value = 2445.2142
$ws.Range("M134").Value = 89.78579999999965
$ws.Range("H141").Value = 31376.434
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 31376.434
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 31376.434
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -41736.434

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 536.6667
$ws.Range("J9").Value = 536.6667
$ws.Range("L9").Value = 1610.0001
$ws.Range("N9").Value = -2058.0001
$ws.Range("H22").Value = 17283.334
$ws.Range("I22").Value = 50550
$ws.Range("J22").Value = 650
$ws.Range("K22").Value = 151650
$ws.Range("L22").Value = 1950
$ws.Range("M22").Value = -151481
$ws.Range("N22").Value = -2288
$ws.Range("H27").Value = 17283.334
$ws.Range("I27").Value = 50550
$ws.Range("J27").Value = 650
$ws.Range("K27").Value = 151650
$ws.Range("L27").Value = 1950
$ws.Range("M27").Value = -151548
$ws.Range("N27").Value = -2154
$ws.Range("H32").Value = 1583.0834
$ws.Range("I32").Value = 400
$ws.Range("J32").Value = 1690.6364
$ws.Range("K32").Value = 1200
$ws.Range("L32").Value = 5071.9092
$ws.Range("M32").Value = -917
$ws.Range("N32").Value = -5637.9092
$ws.Range("H34").Value = 686.1875
$ws.Range("J34").Value = 764.9167
$ws.Range("L34").Value = 2294.7501
$ws.Range("N34").Value = -2462.7501
$ws.Range("H39").Value = 3234.65
$ws.Range("J39").Value = 3234.65
$ws.Range("L39").Value = 9703.950000000001
$ws.Range("N39").Value = -10291.95
$ws.Range("H40").Value = 107.833336
$ws.Range("I40").Value = 116
$ws.Range("J40").Value = 99.666664
$ws.Range("K40").Value = 464
$ws.Range("L40").Value = 398.666656
$ws.Range("M40").Value = -395
$ws.Range("N40").Value = -536.666656
$ws.Range("H55").Value = 4500
$ws.Range("J55").Value = 4500
$ws.Range("L55").Value = 13500
$ws.Range("N55").Value = -13854
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").ClearContents()
$ws.Range("H92").Value = 31250850
$ws.Range("I92").Value = 62500200
$ws.Range("K92").Value = 187500600
$ws.Range("M92").Value = -187499352
$ws.Range("H98").Value = 1883.3334
$ws.Range("I98").Value = 1166.6666
$ws.Range("J98").Value = 2600
$ws.Range("K98").Value = 3499.9998
$ws.Range("L98").Value = 7800
$ws.Range("M98").Value = -2001.9998
$ws.Range("N98").Value = -10796
$ws.Range("H104").Value = 3125
$ws.Range("I104").Value = 2142.8572
$ws.Range("K104").Value = 6428.571599999999
$ws.Range("M104").Value = -3807.571599999999
$ws.Range("H106").Value = 3484.5454
$ws.Range("J106").Value = 3484.5454
$ws.Range("L106").Value = 10453.6362
$ws.Range("N106").Value = -12345.6362
$ws.Range("H107").Value = 7088.3213
$ws.Range("J107").Value = 538.2
$ws.Range("L107").Value = 1614.6
$ws.Range("N107").Value = -5454.6
$ws.Range("H113").Value = 902.1429000000001
$ws.Range("H131").Value = 644.48
$ws.Range("I131").Value = 309.92593
$ws.Range("J131").Value = 768.2192
$ws.Range("K131").Value = 929.77779
$ws.Range("L131").Value = 2304.6576
$ws.Range("M131").Value = 4110.22221
$ws.Range("N131").Value = -12384.6576

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 203
$ws.Range("I3").Value = 203
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 203
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -87
$ws.Range("N3").ClearContents()
$ws.Range("H132").Value = 15588.904
$ws.Range("I132").Value = 4623.6
$ws.Range("J132").Value = 25557.363
$ws.Range("K132").Value = 13870.8
$ws.Range("L132").Value = 76672.08900000001
$ws.Range("M132").Value = -11340.8
$ws.Range("N132").Value = -81732.08900000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 156.66667
$ws.Range("I55").Value = 127.875
$ws.Range("J55").Value = 179.7
$ws.Range("K55").Value = 127.875
$ws.Range("L55").Value = 179.7
$ws.Range("M55").Value = 45.125
$ws.Range("N55").Value = -525.7
$ws.Range("H61").Value = 3294.9092
$ws.Range("I61").Value = 1746.4736
$ws.Range("J61").Value = 13101.667
$ws.Range("K61").Value = 1746.4736
$ws.Range("L61").Value = 13101.667
$ws.Range("M61").Value = -1544.4736
$ws.Range("N61").Value = -13505.667
$ws.Range("H113").Value = 3294.9092
$ws.Range("I113").Value = 1746.4736
$ws.Range("J113").Value = 13101.667
$ws.Range("K113").Value = 1746.4736
$ws.Range("L113").Value = 13101.667
$ws.Range("M113").Value = 423.5264
$ws.Range("N113").Value = -17441.667

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1096.6666
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 1096.6666
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H132").Value = 625.4
$ws.Range("I132").Value = 504.53125
$ws.Range("K132").Value = 1513.59375
$ws.Range("M132").Value = 1016.40625
